$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (shifts existing rows 14-99 down to 15-100)
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new weekly data entry
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 45149
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112009
$ws.Range("G14").Value = "Acelga"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 370
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 941
$ws.Range("N14").Value = "$/atado"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 941
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
